$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry holds the new "Price" (D) and "Volume(1h)" (E) text for one coin
# row. D = $null means the Price column is unchanged for that row (only the
# volume percentage moved). Values are written as literal text (matching the
# source data's inlineStr cells) rather than numbers/percentages.
$updates = @(
    @{ Row = 2; D = "61.489.42"; E = "  +1.10%  " },
    @{ Row = 3; D = "2.382.23"; E = "  +1.14%  " },
    @{ Row = 4; D = $null; E = "  -0.05%  " },
    @{ Row = 5; D = "553.41"; E = "  +2.63%  " },
    @{ Row = 6; D = "140.04"; E = "  +2.01%  " },
    @{ Row = 7; D = $null; E = "  -0.06%  " },
    @{ Row = 8; D = $null; E = "  +0.91%  " },
    @{ Row = 9; D = "2.382.63"; E = "  +1.26%  " },
    @{ Row = 10; D = $null; E = "  +3.82%  " },
    @{ Row = 11; D = $null; E = "  +2.25%  " },
    @{ Row = 12; D = $null; E = "  +2.56%  " },
    @{ Row = 13; D = "0.353"; E = "  +3.60%  " },
    @{ Row = 14; D = "25.65"; E = "  +3.57%  " },
    @{ Row = 15; D = $null; E = "  +6.97%  " },
    @{ Row = 16; D = "2.808.05"; E = "  +0.96%  " },
    @{ Row = 17; D = "61.398.80"; E = "  +1.15%  " },
    @{ Row = 18; D = "2.380.44"; E = "  +1.00%  " },
    @{ Row = 19; D = "10.97"; E = "  +3.79%  " },
    @{ Row = 20; D = $null; E = "  +2.92%  " },
    @{ Row = 21; D = "320.99"; E = "  +1.77%  " },
    @{ Row = 22; D = "6.70"; E = "  +2.09%  " },
    @{ Row = 23; D = $null; E = "  +0.37%  " },
    @{ Row = 24; D = "64.28"; E = "  +1.87%  " },
    @{ Row = 25; D = $null; E = "  -7.78%  " },
    @{ Row = 26; D = "8.90"; E = "  +5.71%  " },
    @{ Row = 27; D = "0.998"; E = "  -0.22%  " },
    @{ Row = 28; D = "2.498.34"; E = "  +1.09%  " },
    @{ Row = 29; D = "8.21"; E = "  +3.51%  " },
    @{ Row = 30; D = "524.49"; E = "  +4.63%  " },
    @{ Row = 31; D = "0.0₃0906"; E = "  +1.61%  " },
    @{ Row = 32; D = $null; E = "  +1.39%  " },
    @{ Row = 33; D = "0.149"; E = "  +2.78%  " },
    @{ Row = 34; D = $null; E = "  +3.41%  " },
    @{ Row = 35; D = "1.52"; E = "  -0.42%  " },
    @{ Row = 36; D = $null; E = "  +0.00%  " },
    @{ Row = 37; D = "5.57"; E = "  +6.14%  " },
    @{ Row = 38; D = "4.71"; E = "  +3.37%  " },
    @{ Row = 39; D = $null; E = "  +6.02%  " },
    @{ Row = 40; D = "0.379"; E = "  +1.98%  " },
    @{ Row = 41; D = "18.55"; E = "  +0.23%  " },
    @{ Row = 42; D = "145.76"; E = "  +5.18%  " },
    @{ Row = 43; D = $null; E = "  -0.01%  " },
    @{ Row = 44; D = "41.40"; E = "  +3.17%  " },
    @{ Row = 45; D = "147.66"; E = "  +6.52%  " },
    @{ Row = 46; D = "2.17"; E = "  +3.70%  " },
    @{ Row = 47; D = "3.61"; E = "  +3.08%  " },
    @{ Row = 48; D = "0.0525"; E = "  +3.09%  " },
    @{ Row = 49; D = "19.90"; E = "  +2.17%  " },
    @{ Row = 50; D = "0.583"; E = "  +2.71%  " },
    @{ Row = 51; D = "0.0909"; E = "  +1.50%  " }
)

foreach ($u in $updates) {
    $r = $u.Row

    if ($null -ne $u.D) {
        $dCell = $ws.Range("D$r")
        $dValue = $u.D
        # Many prices are plain numeric-looking strings (e.g. "553.41").
        # Assigning those straight to .Value would make Excel reinterpret
        # them as numbers, so quote-prefix forces text entry, then the
        # style is reset to Normal so no stray number format lingers.
        $looksNumeric = $dValue -match '^[0-9]+(\.[0-9]+)?$'
        if ($looksNumeric) {
            $dCell.Value = "'" + $dValue
            $dCell.Style = "Normal"
        } else {
            $dCell.Value = $dValue
        }
    }

    $ws.Range("E$r").Value = $u.E
}
